$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 8372.5
$ws.Range("I31").Value = 8372.5
$ws.Range("K31").Value = 25117.5
$ws.Range("M31").Value = -24887.5
$ws.Range("H80").Value = 7187.8125
$ws.Range("I80").Value = 308.44446
$ws.Range("K80").Value = 925.33338
$ws.Range("M80").Value = 72.66661999999997
$ws.Range("H83").Value = 7187.8125
$ws.Range("I83").Value = 308.44446
$ws.Range("K83").Value = 2776.00014
$ws.Range("M83").Value = 2215.99986
$ws.Range("H127").Value = 1296.3077
$ws.Range("I127").Value = 396.75
$ws.Range("J127").Value = 2735.6
$ws.Range("K127").Value = 1190.25
$ws.Range("L127").Value = 8206.8
$ws.Range("M127").Value = 3769.75
$ws.Range("N127").Value = -18126.8
$ws.Range("H129").Value = 842.59
$ws.Range("I129").Value = 349.25
$ws.Range("J129").Value = 936.5595
$ws.Range("K129").Value = 1047.75
$ws.Range("L129").Value = 2809.6785
$ws.Range("M129").Value = 3952.25
$ws.Range("N129").Value = -12809.6785
$ws.Range("H138").Value = 3740.309
$ws.Range("J138").Value = 5138.853
$ws.Range("L138").Value = 15416.559
$ws.Range("N138").Value = -25696.559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2670.318
$ws.Range("I61").Value = 2279.2942
$ws.Range("J61").Value = 3999.8
$ws.Range("K61").Value = 2279.2942
$ws.Range("L61").Value = 3999.8
$ws.Range("M61").Value = -2067.2942
$ws.Range("N61").Value = -4423.8
$ws.Range("H63").Value = 4736.273
$ws.Range("I63").Value = 3699.8333
$ws.Range("K63").Value = 3699.8333
$ws.Range("M63").Value = -3013.8333
$ws.Range("H66").Value = 4736.273
$ws.Range("I66").Value = 3699.8333
$ws.Range("K66").Value = 18499.1665
$ws.Range("M66").Value = -15067.1665
$ws.Range("H88").Value = 2594.1538
$ws.Range("I88").Value = 1963.4
$ws.Range("J88").Value = 2988.375
$ws.Range("K88").Value = 1963.4
$ws.Range("L88").Value = 2988.375
$ws.Range("M88").Value = -1557.4
$ws.Range("N88").Value = -3800.375
$ws.Range("H91").Value = 2594.1538
$ws.Range("I91").Value = 1963.4
$ws.Range("J91").Value = 2988.375
$ws.Range("K91").Value = 1963.4
$ws.Range("L91").Value = 2988.375
$ws.Range("M91").Value = -559.4000000000001
$ws.Range("N91").Value = -5796.375
$ws.Range("H122").Value = 10761.385
$ws.Range("I122").Value = 10836.32
$ws.Range("J122").Value = 8888
$ws.Range("K122").Value = 32508.96
$ws.Range("L122").Value = 26664
$ws.Range("M122").Value = -30058.96
$ws.Range("N122").Value = -31564
$ws.Range("H132").Value = 9789.4375
$ws.Range("I132").Value = 13063.3
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 39189.89999999999
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -36659.89999999999
$ws.Range("N132").Value = -18059
$ws.Range("H136").Value = 2670.318
$ws.Range("I136").Value = 2279.2942
$ws.Range("J136").Value = 3999.8
$ws.Range("K136").Value = 6837.882599999999
$ws.Range("L136").Value = 11999.4
$ws.Range("M136").Value = -4287.882599999999
$ws.Range("N136").Value = -17099.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 32986.41
$ws.Range("I107").Value = 58566.555
$ws.Range("J107").Value = 4208.75
$ws.Range("K107").Value = 58566.555
$ws.Range("L107").Value = 4208.75
$ws.Range("M107").Value = -56646.555
$ws.Range("N107").Value = -8048.75
$ws.Range("H134").Value = 2949.516
$ws.Range("I134").Value = 2378.0952
$ws.Range("K134").Value = 7134.285600000001
$ws.Range("M134").Value = -4599.285600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 50000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""
$ws.Range("H31").Value = 5533.8823
$ws.Range("I31").Value = 2714
$ws.Range("J31").Value = 8040.4443
$ws.Range("K31").Value = 2714
$ws.Range("L31").Value = 8040.4443
$ws.Range("M31").Value = -2419
$ws.Range("N31").Value = -8630.4443
$ws.Range("H34").Value = 5533.8823
$ws.Range("I34").Value = 2714
$ws.Range("J34").Value = 8040.4443
$ws.Range("K34").Value = 2714
$ws.Range("L34").Value = 8040.4443
$ws.Range("M34").Value = -2512
$ws.Range("N34").Value = -8444.4443
$ws.Range("H58").Value = 2016
$ws.Range("I58").Value = 1748
$ws.Range("J58").Value = 2619
$ws.Range("K58").Value = 1748
$ws.Range("L58").Value = 2619
$ws.Range("M58").Value = -1545
$ws.Range("N58").Value = -3025
$ws.Range("H107").Value = 384.2
$ws.Range("I107").Value = 224.4
$ws.Range("J107").Value = 544
$ws.Range("K107").Value = 224.4
$ws.Range("L107").Value = 544
$ws.Range("M107").Value = 1695.6
$ws.Range("N107").Value = -4384
$ws.Range("H123").Value = 37424
$ws.Range("J123").Value = 37424
$ws.Range("L123").Value = 37424
$ws.Range("N123").Value = -47224
$ws.Range("H132").Value = 2205.0833
$ws.Range("I132").Value = 1601.5
$ws.Range("J132").Value = 3412.25
$ws.Range("K132").Value = 4804.5
$ws.Range("L132").Value = 10236.75
$ws.Range("M132").Value = -2274.5
$ws.Range("N132").Value = -15296.75
$ws.Range("H136").Value = 2016
$ws.Range("I136").Value = 1748
$ws.Range("J136").Value = 2619
$ws.Range("K136").Value = 5244
$ws.Range("L136").Value = 7857
$ws.Range("M136").Value = -2694
$ws.Range("N136").Value = -12957

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 836.381
$ws.Range("I107").Value = 1047.8
$ws.Range("J107").Value = 770.3125
$ws.Range("K107").Value = 3143.4
$ws.Range("L107").Value = 2310.9375
$ws.Range("M107").Value = -1223.4
$ws.Range("N107").Value = -6150.9375
$ws.Range("H122").Value = 1112.7333
$ws.Range("J122").Value = 1232.375
$ws.Range("L122").Value = 11091.375
$ws.Range("N122").Value = -15991.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 28991.889
$ws.Range("I97").Value = 42868.418
$ws.Range("J97").Value = 1238.8334
$ws.Range("K97").Value = 42868.418
$ws.Range("L97").Value = 1238.8334
$ws.Range("M97").Value = -42372.418
$ws.Range("N97").Value = -2230.8334
$ws.Range("H102").Value = 3739.7666
$ws.Range("I102").Value = 3603.577
$ws.Range("J102").Value = 4625
$ws.Range("K102").Value = 3603.577
$ws.Range("L102").Value = 4625
$ws.Range("M102").Value = -1981.577
$ws.Range("N102").Value = -7869
$ws.Range("H132").Value = 3951.36
$ws.Range("I132").Value = 3072.6667
$ws.Range("J132").Value = 4228.8423
$ws.Range("K132").Value = 9218.000100000001
$ws.Range("L132").Value = 12686.5269
$ws.Range("M132").Value = -6688.000100000001
$ws.Range("N132").Value = -17746.5269

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3719.1562
$ws.Range("I132").Value = 3748.261
$ws.Range("J132").Value = 3644.7778
$ws.Range("K132").Value = 11244.783
$ws.Range("L132").Value = 10934.3334
$ws.Range("M132").Value = -8714.783
$ws.Range("N132").Value = -15994.3334
$ws.Range("H136").Value = 2990.2
$ws.Range("I136").Value = 1785.3
$ws.Range("K136").Value = 5355.9
$ws.Range("M136").Value = -2805.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 2340
$ws.Range("I29").Value = 2340
$ws.Range("K29").Value = 2340
$ws.Range("M29").Value = -2050
$ws.Range("H81").Value = 78005
$ws.Range("I81").Value = 161825
$ws.Range("J81").Value = 4662.5
$ws.Range("K81").Value = 323650
$ws.Range("L81").Value = 9325
$ws.Range("M81").Value = -322589
$ws.Range("N81").Value = -11447
$ws.Range("H84").Value = 78005
$ws.Range("I84").Value = 161825
$ws.Range("J84").Value = 4662.5
$ws.Range("K84").Value = 1618250
$ws.Range("L84").Value = 46625
$ws.Range("M84").Value = -1612946
$ws.Range("N84").Value = -57233
$ws.Range("H136").Value = 2331.6
$ws.Range("I136").Value = 2059.0435
$ws.Range("J136").Value = 2854
$ws.Range("K136").Value = 6177.130500000001
$ws.Range("L136").Value = 8562
$ws.Range("M136").Value = -3627.130500000001
$ws.Range("N136").Value = -13662
